$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the FROM / TO dates used to compute the accrual period ---
$ws.Range("E15").Formula = "=DATE(2020,11,1)"
$ws.Range("E16").Formula = "=DATE(2020,11,25)"

# --- Clean out the now-obsolete reconciliation rows (old E18, E19, E22) ---
$ws.Range("E18").Clear()
$ws.Range("E19").Clear()
$ws.Range("E22").Clear()

# --- Re-enter the interest total in its new home (E19), same Comma style
#     used elsewhere in the sheet (e.g. E13) ---
$ws.Range("E19").Formula = "=E17*E12"
$ws.Range("E19").NumberFormat = $ws.Range("E13").NumberFormat

# --- Leave a formatted (but empty) placeholder cell at E23 ---
$ws.Range("E23").NumberFormat = $ws.Range("E13").NumberFormat

# --- New helper cell P16, formatted like E13/E19 ---
$ws.Range("P16").NumberFormat = $ws.Range("E13").NumberFormat

# --- Column P widens out now that it holds real content ---
$ws.Columns("P").ColumnWidth = 12.5

# --- Match the recorded UI selection ---
$null = $ws.Range("P16").Select()
